$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 17-19: copy formatting down from the matching template rows first ---
# Row 17 continues row 16's banding (plain Calibri style).
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)

# Rows 18-19 continue row 15's banding (wrapped Cambria style).
$ws.Range("A15:G15").Copy()
$ws.Range("A18:G19").PasteSpecial(-4122)

# Row 19's Notes cell (G19) stays empty and uses the plain style like G16/G17 instead.
$ws.Range("G16").Copy()
$ws.Range("G19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 19 needs the taller 30pt height (matches the other wrapped/long rows).
$ws.Rows.Item(19).RowHeight = 30

# --- Status updates on existing rows ---
$ws.Range("B2").Value = "closed"
$ws.Range("B3").Value = "review"
$ws.Range("B5").Value = "closed"

# --- Row 17 values (new "Controls Team" action item) ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "open"
$ws.Range("C17").Value = 42263
$ws.Range("D17").Value = 42255
$ws.Range("E17").Value = "Controls Team"
$ws.Range("F17").Value = "Develop controller models"
$ws.Range("G17").Value = "Presentation will be made to Dr. Frew on 9-15"

# --- Row 15: "?" placeholder item becomes a real Due date + updated owner/description ---
$ws.Range("C15").Value = 42258
$ws.Range("C15").NumberFormat = "d-mmm"
$ws.Range("E15").Value = "SLAM Team"
$ws.Range("F15").Value = "Trade Study on Hardware (onboard comps)"

# --- Row 18 values (new "SLAM Team" action item) ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "open"
$ws.Range("C18").Value = 42270
$ws.Range("D18").Value = 42255
$ws.Range("E18").Value = "SLAM Team"
$ws.Range("F18").Value = "Prototype SLAM Algorithm on Laptop"
$ws.Range("G18").Value = "Presentation will be made to Dr. Frew on 9-23"

# --- Row 19 values (new "Drew" action item) ---
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "open"
$ws.Range("C19").Value = 42278
$ws.Range("D19").Value = 42255
$ws.Range("E19").Value = "Drew"
$ws.Range("F19").Value = "Order and receive onboard computers and sensors"

# Ensure the date columns keep the d-mmm display used throughout the sheet.
$ws.Range("C17:D19").NumberFormat = "d-mmm"

# Final cursor position left where the author's edit session ended.
$ws.Range("E26").Select()
